{"js": "// Replace each two-digit-by-two-digit multiplication equation in the\n// practice table with its updated version. Each \"from\" text is unique\n// within the document, so a literal (non-wildcard) search locates the\n// single run that needs updating and insertText(..., \"Replace\") swaps\n// its contents in place, preserving the run's formatting (font/size).\nconst replacements = [\n  [\"63\u00d733=2079\", \"83\u00d727=2241\"],\n  [\"37\u00d772=2664\", \"47\u00d719=893\"],\n  [\"39\u00d720=780\", \"80\u00d719=1520\"],\n  [\"73\u00d793=6789\", \"31\u00d746=1426\"],\n  [\"47\u00d772=3384\", \"21\u00d783=1743\"],\n  [\"81\u00d783=6723\", \"49\u00d755=2695\"],\n  [\"98\u00d739=3822\", \"96\u00d745=4320\"],\n  [\"87\u00d738=3306\", \"27\u00d727=729\"],\n  [\"37\u00d718=666\", \"55\u00d797=5335\"],\n  [\"14\u00d714=196\", \"47\u00d788=4136\"],\n  [\"38\u00d729=1102\", \"51\u00d713=663\"],\n  [\"35\u00d742=1470\", \"40\u00d777=3080\"],\n  [\"68\u00d756=3808\", \"25\u00d714=350\"],\n  [\"86\u00d749=4214\", \"17\u00d786=1462\"],\n  [\"24\u00d761=1464\", \"54\u00d743=2322\"],\n  [\"42\u00d768=2856\", \"77\u00d743=3311\"],\n  [\"19\u00d770=1330\", \"63\u00d750=3150\"],\n  [\"93\u00d776=7068\", \"46\u00d762=2852\"],\n  [\"12\u00d756=672\", \"15\u00d778=1170\"],\n  [\"39\u00d777=3003\", \"63\u00d772=4536\"],\n  [\"33\u00d719=627\", \"97\u00d760=5820\"],\n  [\"29\u00d726=754\", \"79\u00d749=3871\"],\n  [\"39\u00d796=3744\", \"95\u00d714=1330\"],\n  [\"55\u00d770=3850\", \"35\u00d731=1085\"],\n  [\"76\u00d772=5472\", \"49\u00d766=3234\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + from);\n  }\n\n  for (const item of results.items) {\n    item.insertText(to, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-by-two-digit multiplication equation in the\n# practice table with its updated version. Content.Find.Execute with\n# Replace:=wdReplaceAll (2) performs a literal (non-wildcard) search and\n# in-place substitution for each unique \"from\" string; run formatting\n# (font/size) is untouched because Word's Find/Replace only swaps text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"63\u00d733=2079\", \"83\u00d727=2241\"),\n    @(\"37\u00d772=2664\", \"47\u00d719=893\"),\n    @(\"39\u00d720=780\", \"80\u00d719=1520\"),\n    @(\"73\u00d793=6789\", \"31\u00d746=1426\"),\n    @(\"47\u00d772=3384\", \"21\u00d783=1743\"),\n    @(\"81\u00d783=6723\", \"49\u00d755=2695\"),\n    @(\"98\u00d739=3822\", \"96\u00d745=4320\"),\n    @(\"87\u00d738=3306\", \"27\u00d727=729\"),\n    @(\"37\u00d718=666\", \"55\u00d797=5335\"),\n    @(\"14\u00d714=196\", \"47\u00d788=4136\"),\n    @(\"38\u00d729=1102\", \"51\u00d713=663\"),\n    @(\"35\u00d742=1470\", \"40\u00d777=3080\"),\n    @(\"68\u00d756=3808\", \"25\u00d714=350\"),\n    @(\"86\u00d749=4214\", \"17\u00d786=1462\"),\n    @(\"24\u00d761=1464\", \"54\u00d743=2322\"),\n    @(\"42\u00d768=2856\", \"77\u00d743=3311\"),\n    @(\"19\u00d770=1330\", \"63\u00d750=3150\"),\n    @(\"93\u00d776=7068\", \"46\u00d762=2852\"),\n    @(\"12\u00d756=672\", \"15\u00d778=1170\"),\n    @(\"39\u00d777=3003\", \"63\u00d772=4536\"),\n    @(\"33\u00d719=627\", \"97\u00d760=5820\"),\n    @(\"29\u00d726=754\", \"79\u00d749=3871\"),\n    @(\"39\u00d796=3744\", \"95\u00d714=1330\"),\n    @(\"55\u00d770=3850\", \"35\u00d731=1085\"),\n    @(\"76\u00d772=5472\", \"49\u00d766=3234\")\n)\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n\n    $found = $rng.Find.Execute(\n        $from,      # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $to,        # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $from\"\n    }\n}\n"}
